$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.970.62"
$ws.Range("D3").Value = "'1.653.62"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'215.05"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'20.22"
$ws.Range("E10").Value = "  +5.10%  "
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("D12").Value = "'1.886.72"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").Value = "'1.647.34"
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "'65.27"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "'26.970.25"
$ws.Range("D18").Value = "'236.63"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.0₃0734"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.77"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +3.67%  "
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").Value = "'145.04"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "'15.85"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").Value = "'1.547.79"
$ws.Range("E32").Value = "  +3.76%  "
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("D34").Value = "'3.09"
$ws.Range("E34").Value = "  +4.98%  "
$ws.Range("E35").Value = "  +10.08%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  +3.41%  "
$ws.Range("E38").Value = "  +9.63%  "
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("D40").Value = "'6.00"
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'66.47"
$ws.Range("E42").Value = "  +9.16%  "
$ws.Range("D43").Value = "'2.23"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").Value = "'1.795.35"
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'0.938"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'0.775"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").Value = "'89.93"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "'1.53"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("D49").Value = "'0.0989"
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "'7.65"
$ws.Range("E51").Value = "  +2.59%  "
